# Update "want to go" counts ($F column) to values captured at gh-pages build 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 6  # was 3
$ws.Cells.Item(4, 6).Value = 513  # was 512
$ws.Cells.Item(5, 6).Value = 13684  # was 13672
$ws.Cells.Item(7, 6).Value = 33  # was 31
$ws.Cells.Item(8, 6).Value = 1754  # was 1751
$ws.Cells.Item(9, 6).Value = 159  # was 157
$ws.Cells.Item(12, 6).Value = 43  # was 42
$ws.Cells.Item(15, 6).Value = 13689  # was 13672
$ws.Cells.Item(16, 6).Value = 349  # was 344
$ws.Cells.Item(18, 6).Value = 9032  # was 9026
$ws.Cells.Item(20, 6).Value = 8152  # was 8139
$ws.Cells.Item(22, 6).Value = 17  # was 16
$ws.Cells.Item(31, 6).Value = 401  # was 400
$ws.Cells.Item(33, 6).Value = 212  # was 211
$ws.Cells.Item(34, 6).Value = 202  # was 200
$ws.Cells.Item(35, 6).Value = 386  # was 385
$ws.Cells.Item(37, 6).Value = 2078  # was 24

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 43  # was 42

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 6  # was 3
$ws.Cells.Item(4, 6).Value = 513  # was 512
$ws.Cells.Item(5, 6).Value = 13684  # was 13672
$ws.Cells.Item(7, 6).Value = 33  # was 31
$ws.Cells.Item(8, 6).Value = 1754  # was 1751
$ws.Cells.Item(9, 6).Value = 159  # was 157
$ws.Cells.Item(12, 6).Value = 43  # was 42
$ws.Cells.Item(15, 6).Value = 13689  # was 13672
$ws.Cells.Item(16, 6).Value = 349  # was 344
$ws.Cells.Item(18, 6).Value = 9032  # was 9026
$ws.Cells.Item(20, 6).Value = 8152  # was 8139
$ws.Cells.Item(22, 6).Value = 17  # was 16
$ws.Cells.Item(31, 6).Value = 43  # was 42
$ws.Cells.Item(33, 6).Value = 401  # was 400
$ws.Cells.Item(35, 6).Value = 212  # was 211
$ws.Cells.Item(36, 6).Value = 202  # was 200
$ws.Cells.Item(37, 6).Value = 386  # was 385
$ws.Cells.Item(39, 6).Value = 2105  # was 24

